# Appends/refreshes the lancers.jp job-listing rows on the "ランサーズ" sheet.
# Mirrors a scrape run at 2026-02-03 01:57:19: one new row inserted after the
# existing 6th entry, five new rows appended at the bottom, and every row's
# "取得日時" (fetched-at) timestamp bumped to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy an existing, already-styled hyperlink cell so new URL cells can pick up
# the exact same "Hyperlink" cell format (Hyperlinks.Add on its own restyles the
# cell with a freshly-minted format record instead of reusing the existing one).
$ws.Range("F2").Copy() | Out-Null

# Row 2
$ws.Cells.Item(2, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(2, 2).Value = 'Gmail、スプレッドシート、Google Driveを連携した 業務効率化システム開発'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5484339'
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5484339') | Out-Null
$ws.Cells.Item(2, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(2, 7).Value = 453
$ws.Cells.Item(2, 8).Value = '🔥AI,Ai ◆効率化,開発'

# Row 3
$ws.Cells.Item(3, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(3, 2).Value = '【完全在宅】AI×Web開発エンジニア募集!業務自動化・AI機能開発'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5483480'
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5483480') | Out-Null
$ws.Cells.Item(3, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 7).Value = 435
$ws.Cells.Item(3, 8).Value = '🔥AI,Ai ◆開発,自動化'

# Row 4
$ws.Cells.Item(4, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(4, 2).Value = '製造業向け設備要件定義書の自動生成AIシステムの開発・DB設計支援エンジニア(AI/バックエンド)'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5473648'
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5473648') | Out-Null
$ws.Cells.Item(4, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 7).Value = 390
$ws.Cells.Item(4, 8).Value = '🔥AI,Ai ◆開発'

# Row 5
$ws.Cells.Item(5, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(5, 2).Value = '※急募:Next.jsによる業務アプリの開発(+Flutter)'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5483966'
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5483966') | Out-Null
$ws.Cells.Item(5, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 7).Value = 225
$ws.Cells.Item(5, 8).Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 6
$ws.Cells.Item(6, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(6, 2).Value = '※急募:Flutterによる業務アプリの開発(+next.js)'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5483967'
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5483967') | Out-Null
$ws.Cells.Item(6, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(6, 7).Value = 218
$ws.Cells.Item(6, 8).Value = '🔥Next.js ◆開発 ◇アプリ'

# Row 7
$ws.Cells.Item(7, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(7, 2).Value = '【急募】楽天RPP広告自動化ツールの開発依頼'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5484001'
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5484001') | Out-Null
$ws.Cells.Item(7, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(7, 7).Value = 213
$ws.Cells.Item(7, 8).Value = '◆ツール,開発'

# Row 8
$ws.Cells.Item(8, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(8, 2).Value = '急募 【急募】日繰資金繰表自動算出ツールの開発依頼'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5484458'
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5484458') | Out-Null
$ws.Cells.Item(8, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 7).Value = 123
$ws.Cells.Item(8, 8).Value = '◆ツール,開発'

# Row 9
$ws.Cells.Item(9, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(9, 2).Value = 'X(旧twitter)のロック解除自動化システム構築'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5484231'
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5484231') | Out-Null
$ws.Cells.Item(9, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(9, 7).Value = 103
$ws.Cells.Item(9, 8).Value = '◆自動化'

# Row 10
$ws.Cells.Item(10, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(10, 2).Value = '(仕様削減)【受注メールを元にしたスクレピング&抽出情報管理ツール】'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5484248'
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5484248') | Out-Null
$ws.Cells.Item(10, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 7).Value = 98
$ws.Cells.Item(10, 8).Value = '◆ツール ◇管理'

# Row 11
$ws.Cells.Item(11, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(11, 2).Value = 'ストレスチェック集団分析自動化・レポート出力システムを探しています'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5484413'
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5484413') | Out-Null
$ws.Cells.Item(11, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(11, 7).Value = 98
$ws.Cells.Item(11, 8).Value = '◆自動化'

# Row 12
$ws.Cells.Item(12, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(12, 2).Value = '製造業DXプロダクト開発のプロダクトマネージャー募集'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5468432'
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5468432') | Out-Null
$ws.Cells.Item(12, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(12, 7).Value = 75
$ws.Cells.Item(12, 8).Value = '◆開発'

# Row 13
$ws.Cells.Item(13, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(13, 2).Value = '【急募】新しいWebサービスの開発パートナーを探しています!'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5483482'
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), 'https://www.lancers.jp/work/detail/5483482') | Out-Null
$ws.Cells.Item(13, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(13, 7).Value = 75
$ws.Cells.Item(13, 8).Value = '◆開発'

# Row 14
$ws.Cells.Item(14, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(14, 2).Value = 'リアルタイム顔変換開発案件 技術判断を整理してくれる方(短期スポット)募集'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5484376'
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), 'https://www.lancers.jp/work/detail/5484376') | Out-Null
$ws.Cells.Item(14, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14, 7).Value = 68
$ws.Cells.Item(14, 8).Value = '◆開発'

# Row 15
$ws.Cells.Item(15, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(15, 2).Value = '四柱推命の命式自動計算プログラム(Web/Excel)の開発依頼'
$ws.Cells.Item(15, 3).Value = 'システム開発'
$ws.Cells.Item(15, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = '期限情報なし'
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5484177'
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), 'https://www.lancers.jp/work/detail/5484177') | Out-Null
$ws.Cells.Item(15, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 7).Value = 68
$ws.Cells.Item(15, 8).Value = '◆開発'

# Row 16
$ws.Cells.Item(16, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(16, 2).Value = '美容皮膚科向け LINE連携型BtoB SaaS(MVP) の開発案件'
$ws.Cells.Item(16, 3).Value = 'システム開発'
$ws.Cells.Item(16, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(16, 5).Value = '期限情報なし'
$ws.Cells.Item(16, 6).Value = 'https://www.lancers.jp/work/detail/5483503'
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), 'https://www.lancers.jp/work/detail/5483503') | Out-Null
$ws.Cells.Item(16, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16, 7).Value = 68
$ws.Cells.Item(16, 8).Value = '◆開発'

# Row 17
$ws.Cells.Item(17, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(17, 2).Value = '【急募】Notion×Slackでのオンライン講座運営システム構築'
$ws.Cells.Item(17, 3).Value = 'システム開発'
$ws.Cells.Item(17, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(17, 5).Value = '期限情報なし'
$ws.Cells.Item(17, 6).Value = 'https://www.lancers.jp/work/detail/5483854'
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), 'https://www.lancers.jp/work/detail/5483854') | Out-Null
$ws.Cells.Item(17, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 7).Value = 28

# Row 18
$ws.Cells.Item(18, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(18, 2).Value = '【急募】Unity/Photonでのメタバースマルチプレイヤー機能実装'
$ws.Cells.Item(18, 3).Value = 'システム開発'
$ws.Cells.Item(18, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(18, 5).Value = '期限情報なし'
$ws.Cells.Item(18, 6).Value = 'https://www.lancers.jp/work/detail/5484436'
$ws.Hyperlinks.Add($ws.Cells.Item(18, 6), 'https://www.lancers.jp/work/detail/5484436') | Out-Null
$ws.Cells.Item(18, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(18, 7).Value = 18

# Row 19
$ws.Cells.Item(19, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(19, 2).Value = '防災プラットフォームの作成'
$ws.Cells.Item(19, 3).Value = 'システム開発'
$ws.Cells.Item(19, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(19, 5).Value = '期限情報なし'
$ws.Cells.Item(19, 6).Value = 'https://www.lancers.jp/work/detail/5484389'
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), 'https://www.lancers.jp/work/detail/5484389') | Out-Null
$ws.Cells.Item(19, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19, 7).Value = 18

# Row 20
$ws.Cells.Item(20, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(20, 2).Value = '限定公開 限定公開の仕事'
$ws.Cells.Item(20, 3).Value = 'システム開発'
$ws.Cells.Item(20, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(20, 5).Value = '期限情報なし'
$ws.Cells.Item(20, 6).Value = 'https://www.lancers.jp/work/detail/5484369'
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), 'https://www.lancers.jp/work/detail/5484369') | Out-Null
$ws.Cells.Item(20, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20, 7).Value = 18

# Row 21
$ws.Cells.Item(21, 1).Value = '2026-02-03 01:57:19'
$ws.Cells.Item(21, 2).Value = '【市場調査】海外向けデジタルサービスの価値評価依頼'
$ws.Cells.Item(21, 3).Value = 'システム開発'
$ws.Cells.Item(21, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(21, 5).Value = '期限情報なし'
$ws.Cells.Item(21, 6).Value = 'https://www.lancers.jp/work/detail/5483504'
$ws.Hyperlinks.Add($ws.Cells.Item(21, 6), 'https://www.lancers.jp/work/detail/5483504') | Out-Null
$ws.Cells.Item(21, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(21, 7).Value = 13

$excel.CutCopyMode = $false

